$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect repulled data / recalculated mean
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = -4
$ws.Range("F7").Value = -2
